$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The payment code has now been generated ("Hasil Generate" was the
# placeholder before generation) -> DISK230200212
$kodePembayaran = "DISK230200212"

# N2 (KODE_PEMBAYARAN column) gets the generated code
$ws.Range("N2").Value = $kodePembayaran

# The approving role changed from "Penyelia Teller" to "Penyelia Settlement"
$ws.Range("I2").Value = "Penyelia Settlement"

# PREPARATION (F2) text references the same generated payment code
$ws.Range("F2").Value = "Username : 32193;`nPassword : bni1234;`nRole : 38 - Penyelia Teller;`nKode Pembayaran : " + $kodePembayaran

# Row 2 wraps onto one additional line now that the PREPARATION text is
# longer, so its auto-fit height grows.
$ws.Rows("2:2").RowHeight = 63.75

# Update the active selection shown when the sheet is reopened.
$ws.Range("J2").Select() | Out-Null
